$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 311, shifting existing rows 311..348 down to 312..349.
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new weekly price-survey record.
$ws.Cells.Item(311, 1).Value = 4
$ws.Cells.Item(311, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(311, 3).Value = "Los Lagos"
$ws.Cells.Item(311, 4).Value = 44918
$ws.Cells.Item(311, 5).Value = 10
$ws.Cells.Item(311, 6).Value = "Fruta"
$ws.Cells.Item(311, 7).Value = 100108
$ws.Cells.Item(311, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(311, 9).Value = 100108005
$ws.Cells.Item(311, 10).Value = "Piña"
$ws.Cells.Item(311, 11).Value = "Caramelo"
$ws.Cells.Item(311, 12).Value = "Segunda"
$ws.Cells.Item(311, 13).Value = 300
$ws.Cells.Item(311, 14).Value = 22000
$ws.Cells.Item(311, 15).Value = 23000
$ws.Cells.Item(311, 16).Value = 22500
$ws.Cells.Item(311, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(311, 18).Value = "Ecuador"
$ws.Cells.Item(311, 19).Value = 1607
$ws.Cells.Item(311, 20).Value = 14
